$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: "q" -> "question", "a" -> "answer"
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Move selection to A2 (matches sheetView selection activeCell="A2" in the diff)
$ws.Range("A2").Select()
